$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to match the collaboration database (UCSF)
$ws.Name = "UCSF"

# Update Thongchai Masilela's email address to include the additional
# curie.fr address, keeping the jiscmail / collaboration database in sync
$ws.Range("E2").Value = "Thongchai.Masilela@ucsf.edu; thongchai.masilela@curie.fr"

# Reflect the active selection moving to the edited cell
$ws.Range("E2").Select()
